$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04904614488327752
$ws.Range("D2").Value = 0.1367748167803811
$ws.Range("E2").Value = 0.1464392491620892
$ws.Range("F2").Value = 2.05016538974661
$ws.Range("G2").Value = 0.002434556470929708
$ws.Range("J2").Value = 0.2009711305352084
$ws.Range("O2").Value = 5.513759881882606
$ws.Range("C3").Value = 0.04352966595151031
$ws.Range("D3").Value = 0.1371026173882441
$ws.Range("E3").Value = 0.1446177847972727
$ws.Range("F3").Value = 1.980028762343863
$ws.Range("G3").Value = 0.002440342032459342
$ws.Range("J3").Value = 0.1953332899241502
$ws.Range("O3").Value = 5.291219219213929
$ws.Range("C4").Value = 0.04016292325236748
$ws.Range("D4").Value = 0.1373433076147101
$ws.Range("E4").Value = 0.1435801722260734
$ws.Range("F4").Value = 1.938296182754314
$ws.Range("G4").Value = 0.002444081946923408
$ws.Range("J4").Value = 0.1920013065019575
$ws.Range("O4").Value = 5.158232048430364
$ws.Range("C5").Value = 0.03879591925563375
$ws.Range("D5").Value = 0.1374513150845225
$ws.Range("E5").Value = 0.1431775837973959
$ws.Range("F5").Value = 1.921622456926428
$ws.Range("G5").Value = 0.002445653326873506
$ws.Range("J5").Value = 0.1906758992079887
$ws.Range("O5").Value = 5.104950414272537
$ws.Range("C6").Value = 0.03856922492485637
$ws.Range("D6").Value = 0.137469849400464
$ws.Range("E6").Value = 0.143111955325697
$ws.Range("F6").Value = 1.918873824376178
$ws.Range("G6").Value = 0.002445917117038179
$ws.Range("J6").Value = 0.1904577674833234
$ws.Range("O6").Value = 5.096157901393042
$ws.Range("C7").Value = 0.04014446743859423
$ws.Range("D7").Value = 0.1373447240423431
$ws.Range("E7").Value = 0.1435746608811499
$ws.Range("F7").Value = 1.938069971089789
$ws.Range("G7").Value = 0.002444102947217419
$ws.Range("J7").Value = 0.1919833006641554
$ws.Range("O7").Value = 5.157509791371922
$ws.Range("C8").Value = 0.04713973221738854
$ws.Range("D8").Value = 0.1368796679017592
$ws.Range("E8").Value = 0.1457944037208492
$ws.Range("F8").Value = 2.02570454850482
$ws.Range("G8").Value = 0.002436512508476527
$ws.Range("J8").Value = 0.1990001632839267
$ws.Range("O8").Value = 5.436265195847284
$ws.Range("C9").Value = 0.06102711453299037
$ws.Range("D9").Value = 0.1362800148175261
$ws.Range("E9").Value = 0.1507915675395495
$ws.Range("F9").Value = 2.20823290739537
$ws.Range("G9").Value = 0.002423108042804387
$ws.Range("J9").Value = 0.213799379725117
$ws.Range("O9").Value = 6.012248861076273
$ws.Range("C10").Value = 0.07134550460989431
$ws.Range("D10").Value = 0.1360293435457436
$ws.Range("E10").Value = 0.1548610410016593
$ws.Range("F10").Value = 2.349018070605268
$ws.Range("G10").Value = 0.002414151305276263
$ws.Range("J10").Value = 0.2253216067526012
$ws.Range("O10").Value = 6.45385907212642
$ws.Range("C11").Value = 0.0760672981778896
$ws.Range("D11").Value = 0.1359564470757633
$ws.Range("E11").Value = 0.1567999944200409
$ws.Range("F11").Value = 2.414554617585054
$ws.Range("G11").Value = 0.00241026789955366
$ws.Range("J11").Value = 0.2307078846773578
$ws.Range("O11").Value = 6.658884712930103
$ws.Range("C12").Value = 0.0778595277758285
$ws.Range("D12").Value = 0.1359347498822103
$ws.Range("E12").Value = 0.1575469301407892
$ws.Range("F12").Value = 2.4395891285522
$ws.Range("G12").Value = 0.002408824648682824
$ws.Range("J12").Value = 0.2327686092152987
$ws.Range("O12").Value = 6.737126306207074
$ws.Range("C13").Value = 0.07747334993777599
$ws.Range("D13").Value = 0.1359391601642841
$ws.Range("E13").Value = 0.1573854981983658
$ws.Range("F13").Value = 2.434187794826073
$ws.Range("G13").Value = 0.002409134266518744
$ws.Range("J13").Value = 0.2323238548569009
$ws.Range("O13").Value = 6.720248635271673
$ws.Range("C14").Value = 0.07621466101760177
$ws.Range("D14").Value = 0.1359545436919163
$ws.Range("E14").Value = 0.1568611903521528
$ws.Range("F14").Value = 2.416609852933561
$ws.Range("G14").Value = 0.002410148615919035
$ws.Range("J14").Value = 0.2308769982900429
$ws.Range("O14").Value = 6.665309566532869
$ws.Range("C15").Value = 0.07544422848580723
$ws.Range("D15").Value = 0.1359647355961826
$ws.Range("E15").Value = 0.156541692499772
$ws.Range("F15").Value = 2.405871225399522
$ws.Range("G15").Value = 0.002410773485899977
$ws.Range("J15").Value = 0.2299935064926046
$ws.Range("O15").Value = 6.631736554161307
$ws.Range("C16").Value = 0.07103750082835347
$ws.Range("D16").Value = 0.1360349345471192
$ws.Range("E16").Value = 0.1547360980360324
$ws.Range("F16").Value = 2.344765360487202
$ws.Range("G16").Value = 0.002414408925671347
$ws.Range("J16").Value = 0.2249725343125277
$ws.Range("O16").Value = 6.440544070253679
$ws.Range("C17").Value = 0.06834138805332657
$ws.Range("D17").Value = 0.1360885291582861
$ws.Range("E17").Value = 0.1536509507887409
$ws.Range("F17").Value = 2.30766303358115
$ws.Range("G17").Value = 0.00241668796849577
$ws.Range("J17").Value = 0.2219295774008998
$ws.Range("O17").Value = 6.324318542011838
$ws.Range("C18").Value = 0.06679327129099022
$ws.Range("D18").Value = 0.1361232276551974
$ws.Range("E18").Value = 0.1530350553665514
$ws.Range("F18").Value = 2.286463095253168
$ws.Range("G18").Value = 0.002418016807068295
$ws.Range("J18").Value = 0.2201929597444803
$ws.Range("O18").Value = 6.257857391018661
$ws.Range("C19").Value = 0.06626954953242148
$ws.Range("D19").Value = 0.1361356412918582
$ws.Range("E19").Value = 0.152827938718044
$ws.Range("F19").Value = 2.27930918533167
$ws.Range("G19").Value = 0.00241846982404647
$ws.Range("J19").Value = 0.219607301197982
$ws.Range("O19").Value = 6.235421298198503
$ws.Range("C20").Value = 0.06862812184921552
$ws.Range("D20").Value = 0.1360824232134874
$ws.Range("E20").Value = 0.1537656118770272
$ws.Range("F20").Value = 2.31159809642412
$ws.Range("G20").Value = 0.002416443499720094
$ws.Range("J20").Value = 0.22225209496861
$ws.Range("O20").Value = 6.336650659906809
$ws.Range("C21").Value = 0.07658425336620667
$ws.Range("D21").Value = 0.135949864921713
$ws.Range("E21").Value = 0.1570148470531407
$ws.Range("F21").Value = 2.421767003689439
$ws.Range("G21").Value = 0.002409849936775932
$ws.Range("J21").Value = 0.2313014016789623
$ws.Range("O21").Value = 6.681430075679373
$ws.Range("C22").Value = 0.08180854445863872
$ws.Range("D22").Value = 0.1358976581553293
$ws.Range("E22").Value = 0.159212451786388
$ws.Range("F22").Value = 2.495036253509966
$ws.Range("G22").Value = 0.002405699776436031
$ws.Range("J22").Value = 0.2373385043004532
$ws.Range("O22").Value = 6.910281055858945
$ws.Range("C23").Value = 0.07901794153349329
$ws.Range("D23").Value = 0.1359223743081444
$ws.Range("E23").Value = 0.1580327467679723
$ws.Range("F23").Value = 2.455814227775562
$ws.Range("G23").Value = 0.002407900290276882
$ws.Range("J23").Value = 0.2341050664453803
$ws.Range("O23").Value = 6.787814350591702
$ws.Range("C24").Value = 0.06849848355155075
$ws.Range("D24").Value = 0.1360851716055436
$ws.Range("E24").Value = 0.1537137487798041
$ws.Range("F24").Value = 2.309818647658062
$ws.Range("G24").Value = 0.002416553966181921
$ws.Range("J24").Value = 0.2221062448667368
$ws.Range("O24").Value = 6.331074193523875
$ws.Range("C25").Value = 0.05725082268395454
$ws.Range("D25").Value = 0.1364088544321582
$ws.Range("E25").Value = 0.1493701888532044
$ws.Range("F25").Value = 2.157693419874022
$ws.Range("G25").Value = 0.002426576956399249
$ws.Range("J25").Value = 0.2096829376188794
$ws.Range("O25").Value = 5.853232542874025
